$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns keep their text formatting so that
# numeric-looking strings (e.g. "65.702.83", "0.497", "1.00") are not
# silently converted into numbers by Excel, matching the original inline-string data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "65.702.83"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "3.581.81"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "603.72"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "137.00"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D7").Value = "3.579.32"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  +5.54%  "
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "4.187.46"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "28.14"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "3.584.98"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "65.754.37"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "10.04"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").Value = "5.88"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").Value = "394.92"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "0.589"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").Value = "3.722.87"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").Value = "1.60"
$ws.Range("E29").Value = "  +26.60%  "
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("D31").Value = "8.57"
$ws.Range("E31").Value = "  +5.08%  "
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "3.582.75"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "24.49"
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.149"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").Value = "  +6.87%  "
$ws.Range("E38").Value = "  +4.97%  "
$ws.Range("D39").Value = "7.05"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").Value = "167.78"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").Value = "0.0833"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("D42").Value = "0.837"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").Value = "27.12"
$ws.Range("E43").Value = "  +4.21%  "
$ws.Range("E44").Value = "  +8.10%  "
$ws.Range("D45").Value = "43.09"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "4.53"
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("D49").Value = "7.00"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").Value = "2.455.00"
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("D51").Value = "0.903"
$ws.Range("E51").Value = "  +9.65%  "
